$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1311.9524
$ws.Range("I98").Value = 927.8125
$ws.Range("J98").Value = 2541.2
$ws.Range("K98").Value = 927.8125
$ws.Range("L98").Value = 2541.2
$ws.Range("M98").Value = 570.1875
$ws.Range("N98").Value = -5537.2
$ws.Range("H122").Value = 1311.9524
$ws.Range("I122").Value = 927.8125
$ws.Range("J122").Value = 2541.2
$ws.Range("K122").Value = 2783.4375
$ws.Range("L122").Value = 7623.599999999999
$ws.Range("M122").Value = -333.4375
$ws.Range("N122").Value = -12523.6
$ws.Range("H137").Value = 4879282.5
$ws.Range("I137").Value = 829.38464
$ws.Range("J137").Value = 13335268
$ws.Range("K137").Value = 2488.15392
$ws.Range("L137").Value = 40005804
$ws.Range("M137").Value = 61.84608000000026
$ws.Range("N137").Value = -40010904
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 71435150
$ws.Range("I61").Value = 100006000
$ws.Range("K61").Value = 100006000
$ws.Range("M61").Value = -100005788
$ws.Range("H136").Value = 71435150
$ws.Range("I136").Value = 100006000
$ws.Range("K136").Value = 300018000
$ws.Range("M136").Value = -300015450
$ws.Range("H138").Value = 58985.2
$ws.Range("J138").Value = 58985.2
$ws.Range("L138").Value = 58985.2
$ws.Range("N138").Value = -69265.2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2623.46
$ws.Range("I134").Value = 1505.7
$ws.Range("J134").Value = 4300.1
$ws.Range("K134").Value = 4517.1
$ws.Range("L134").Value = 12900.3
$ws.Range("M134").Value = -1982.1
$ws.Range("N134").Value = -17970.3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9263584
$ws.Range("I31").Value = 4747.517
$ws.Range("J31").Value = 47621620
$ws.Range("K31").Value = 4747.517
$ws.Range("L31").Value = 47621620
$ws.Range("M31").Value = -4452.517
$ws.Range("N31").Value = -47622210
$ws.Range("H34").Value = 9263584
$ws.Range("I34").Value = 4747.517
$ws.Range("J34").Value = 47621620
$ws.Range("K34").Value = 4747.517
$ws.Range("L34").Value = 47621620
$ws.Range("M34").Value = -4545.517
$ws.Range("N34").Value = -47622024
$ws.Range("H58").Value = 2587.4285
$ws.Range("I58").Value = 662.3
$ws.Range("J58").Value = 7400.25
$ws.Range("K58").Value = 662.3
$ws.Range("L58").Value = 7400.25
$ws.Range("M58").Value = -459.3
$ws.Range("N58").Value = -7806.25
$ws.Range("H132").Value = 2705.6562
$ws.Range("I132").Value = 1925.7307
$ws.Range("K132").Value = 5777.1921
$ws.Range("M132").Value = -3247.1921
$ws.Range("H136").Value = 2587.4285
$ws.Range("I136").Value = 662.3
$ws.Range("J136").Value = 7400.25
$ws.Range("K136").Value = 1986.9
$ws.Range("L136").Value = 22200.75
$ws.Range("M136").Value = 563.1000000000001
$ws.Range("N136").Value = -27300.75
$ws.Range("H140").Value = 35105.332
$ws.Range("J140").Value = 35105.332
$ws.Range("L140").Value = 35105.332
$ws.Range("N140").Value = -45465.332
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1380.4667
$ws.Range("I5").Value = 1004
$ws.Range("J5").Value = 1407.3572
$ws.Range("K5").Value = 3012
$ws.Range("L5").Value = 4222.071599999999
$ws.Range("M5").Value = -2900
$ws.Range("N5").Value = -4446.071599999999
$ws.Range("H122").Value = 952.69446
$ws.Range("I122").Value = 798.88
$ws.Range("J122").Value = 1302.2727
$ws.Range("K122").Value = 7189.92
$ws.Range("L122").Value = 11720.4543
$ws.Range("M122").Value = -4739.92
$ws.Range("N122").Value = -16620.4543
$ws.Range("H135").Value = 1380.4667
$ws.Range("I135").Value = 1004
$ws.Range("J135").Value = 1407.3572
$ws.Range("K135").Value = 9036
$ws.Range("L135").Value = 12666.2148
$ws.Range("M135").Value = -6501
$ws.Range("N135").Value = -17736.2148
$ws.Range("H137").Value = 3971688
$ws.Range("I137").Value = 11906387
$ws.Range("J137").Value = 4338.7144
$ws.Range("K137").Value = 35719161
$ws.Range("L137").Value = 13016.1432
$ws.Range("M137").Value = -35714061
$ws.Range("N137").Value = -23216.1432
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2179.1
$ws.Range("I107").Value = 2587.25
$ws.Range("J107").Value = 546.5
$ws.Range("K107").Value = 2587.25
$ws.Range("L107").Value = 546.5
$ws.Range("M107").Value = -667.25
$ws.Range("N107").Value = -4386.5
$ws.Range("H122").Value = 4445912.5
$ws.Range("I122").Value = 5556790
$ws.Range("J122").Value = 2402.6667
$ws.Range("K122").Value = 16670370
$ws.Range("L122").Value = 7208.000100000001
$ws.Range("M122").Value = -16667920
$ws.Range("N122").Value = -12108.0001
$ws.Range("H126").Value = 3353.4707
$ws.Range("I126").Value = 2082.2222
$ws.Range("K126").Value = 6246.6666
$ws.Range("M126").Value = -3776.6666
$ws.Range("H132").Value = 5943.684
$ws.Range("I132").Value = 5472.5
$ws.Range("J132").Value = 6161.154
$ws.Range("K132").Value = 16417.5
$ws.Range("L132").Value = 18483.462
$ws.Range("M132").Value = -13887.5
$ws.Range("N132").Value = -23543.462
$ws.Range("H134").Value = 38775.332
$ws.Range("J134").Value = 38775.332
$ws.Range("L134").Value = 116325.996
$ws.Range("N134").Value = -121395.996
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3025.093
$ws.Range("I40").Value = 3585
$ws.Range("J40").Value = 2438.524
$ws.Range("K40").Value = 3585
$ws.Range("L40").Value = 2438.524
$ws.Range("M40").Value = -3449
$ws.Range("N40").Value = -2710.524
$ws.Range("H132").Value = 11636604
$ws.Range("I132").Value = 6122.64
$ws.Range("J132").Value = 27790050
$ws.Range("K132").Value = 18367.92
$ws.Range("L132").Value = 83370150
$ws.Range("M132").Value = -15837.92
$ws.Range("N132").Value = -83375210
$ws.Range("H139").Value = 41379.5
$ws.Range("J139").Value = 41445.816
$ws.Range("L139").Value = 41445.816
$ws.Range("N139").Value = -51725.816
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1611.0526
$ws.Range("I126").Value = 1114.0541
$ws.Range("J126").Value = 20000
$ws.Range("K126").Value = 3342.1623
$ws.Range("L126").Value = 60000
$ws.Range("M126").Value = -872.1623
$ws.Range("N126").Value = -64940
$ws.Range("H136").Value = 1139.3
$ws.Range("I136").Value = 1099.2222
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 3297.6666
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -747.6665999999996
$ws.Range("N136").Value = -9600
